# "add ID button works" - clicking the sheet's "Add ID" button generates a
# new random e-mail/identity string and writes the newest one into the
# Email cell (E2) of the register sheet. Replicate that visible effect:
# E2's displayed text moves on to the latest generated address while the
# underlying mailto: hyperlink (set up separately) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "HRzSWz3NBH@gmail.com"
